$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 (App server rule): bump the sequence number and describe the rule.
$ws.Range("A3").Value = 150
# Row 4 (DB server rule): bump the sequence number and describe the rule.
$ws.Range("A4").Value = 200

# Description column (J) — written before the other shared strings so the
# rebuilt shared-string table keeps the same ordering the workbook shipped
# with (description text first, then ports, then direction).
$ws.Range("J3").Value = "App server outbound"
$ws.Range("J4").Value = "DB server Outbound"

# TCP Port column (E) switches from a bare numeric "80443" to the text "80, 443".
$ws.Range("E3").Value = "80, 443"
$ws.Range("E4").Value = "80, 443"

# Direction column (I) switches from "outbound" to "Outbound".
$ws.Range("I3").Value = "Outbound"
$ws.Range("I4").Value = "Outbound"

# Restore the cursor/selection to match where the author left it.
$ws.Range("F5").Select()
